$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column A for bias_id; this shifts the existing
# bias / experiment / content / variables columns one to the right
# (B / C / D / E) and keeps their bestFit column-width metadata intact.
$ws.Range("A1").EntireColumn.Insert()

# --- Write new / changed cell values in the same order the strings were
# --- first introduced so the shared-strings table lines up with the
# --- authored edit (first use of a brand new string decides its slot).

# 1) E3 introduces the second JSON "variables" payload.
$ws.Range("E3").Value = '{"this": 841758, "that": 1341}'

# 2) A1 introduces the new "bias_id" header.
$ws.Range("A1").Value = "bias_id"

# 3) C1 replaces "experiment" with "experiment_type".
$ws.Range("C1").Value = "experiment_type"

# 4) D2 introduces the shared sentence used for every data row's content.
$ws.Range("D2").Value = "Let's test this and whether it is different to that."
$ws.Range("D3").Value = "Let's test this and whether it is different to that."

# 5) C4 introduces the new "test" experiment_type value (also a new row).
$ws.Range("C4").Value = "test"

# Remaining (already-known) values.
$ws.Range("A2").Value = 10
$ws.Range("A3").Value = 11
$ws.Range("A4").Value = 12

$ws.Range("B2").Value = "category_size_bias"
$ws.Range("B3").Value = "category_size_bias"
$ws.Range("B4").Value = "category_size_bias"

$ws.Range("C2").Value = "standard"
$ws.Range("C3").Value = "odd_numbers"

$ws.Range("D4").Value = "Let's test this and whether it is different to that."

# E2 keeps the original JSON payload (value unchanged by the edit).
$ws.Range("E2").Value = '{"this": 12, "that": 22}'

# A2 carries an explicit "General" number format in the authored file.
$ws.Range("A2").NumberFormat = "General"

# Column widths for the (re)introduced columns C/D/E; column B keeps the
# width it inherited from the insert above.
$ws.Columns.Item(3).ColumnWidth = 13.833333333333332
$ws.Columns.Item(4).ColumnWidth = 38.5
$ws.Columns.Item(5).ColumnWidth = 24.166666666666668

$ws.Range("E4").Select()
